$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the Instagram profile selector / typing data
$ws.Range("A2").Value = "https://www.instagram.com/oldrich.hanak.3"
$ws.Range("B2").Value = "hex ️"
